$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.539.56"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.216.07"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'240.87"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'75.50"
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'41.37"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'0.0927"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "'55.00"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "2.542.63"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "2.214.94"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "'0.799"
$ws.Range("E18").Value = "  -4.50%  "
$ws.Range("D19").Value = "42.370.18"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "'0.0000105"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'70.74"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").Value = "'5.92"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "'10.04"
$ws.Range("E23").Value = "  -9.27%  "
$ws.Range("D24").Value = "'228.78"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'10.92"
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "'3.36"
$ws.Range("E28").Value = "  -7.14%  "
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").Value = "'2.18"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'172.76"
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("D32").Value = "'34.39"
$ws.Range("E32").Value = "  +13.17%  "
$ws.Range("D33").Value = "'20.24"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "'5.35"
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("D39").Value = "'0.0322"
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("D40").Value = "'12.61"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").Value = "'5.49"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "'60.64"
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").Value = "'0.198"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D46").Value = "'0.0981"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").Value = "'99.59"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.32"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.10"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "'0.424"
$ws.Range("E51").Value = "  +14.75%  "
